# Insert a new data row for Maracuyá (Arica) at row 64, pushing the
# existing rows 64-144 down to 65-145.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value  = 1
$ws.Cells.Item(64, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value  = 44799
$ws.Cells.Item(64, 5).Value  = 15
$ws.Cells.Item(64, 6).Value  = "Fruta"
$ws.Cells.Item(64, 7).Value  = 100108
$ws.Cells.Item(64, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(64, 9).Value  = 100108003
$ws.Cells.Item(64, 10).Value = "Maracuyá"
$ws.Cells.Item(64, 11).Value = "Sin especificar"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 160
$ws.Cells.Item(64, 14).Value = 22000
$ws.Cells.Item(64, 15).Value = 23000
$ws.Cells.Item(64, 16).Value = 22500
$ws.Cells.Item(64, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(64, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 19).Value = 1125
$ws.Cells.Item(64, 20).Value = 20
